$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title placeholder: merge "Testing" / " " / "custom" / " " / "properties" runs
# into a single run by rewriting the full character range.
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Characters(1, $titleRange.Length).Text = "Testing custom properties"

# Subtitle placeholder: keep the two leading line breaks untouched, but merge
# the trailing "A." / " " / "M." runs into a single run.
$subRange = $s.Shapes.Item(2).TextFrame.TextRange
$subRange.Characters(3, $subRange.Length - 2).Text = "A. M."
